$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.02629232850335
$ws.Range("C2").Value = 15.83682600140548
$ws.Range("D2").Value = 16.58386968409785

$ws.Range("B3").Value = 1.253417026148696
$ws.Range("C3").Value = 1.770451996115586
$ws.Range("D3").Value = 2.205118070226383

$ws.Range("B4").Value = 0.254187766061768
$ws.Range("C4").Value = 0.3498719510603343
$ws.Range("D4").Value = 0.4369279524342684

$ws.Range("B5").Value = 81.03109050255415
$ws.Range("C5").Value = 81.9984482992256
$ws.Range("D5").Value = 83.15410987302802
